# article 48 is live
# Adds row 49 to Sheet1: Ser=48, Date=43878 (2020-02-17), Content (blog post),
# Author="Qasim Ali", Tags.
#
# The big HTML-ish blog body is transported base64-encoded to avoid any
# quoting / escaping issues with the PowerShell-style interpreter.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$bodyB64 = "aDE6IFN0b3AgYmVpbmcganVkZ2VtZW50YWwKcC5ub3RlOiBJIGFtIG5lcnZvdXMsIGJlY2F1c2UgaXQgaXMgYSBib2xkIHN0ZXAuIERvbuKAmXQga25vdyB3aGljaCBzaWRlIHRoaXMgYmxvZyBnb2VzLiBsZXRzIHNlZS4KaWZyYW1lLndpZHRoLW0tNzUuZC1mbGV4Lm0tYXV0bzogaHR0cHM6Ly93d3cueW91dHViZS5jb20vZW1iZWQveHRFS0tJMFdtTXMKcC5jYXB0aW9uLnRleHQtY2VudGVyOiBUaGlzIGlzIG1lIHdyaXRpbmcgdGhpcyBibG9nIGluIHRoZSBtb3JuaW5nIGJlZm9yZSBmYWphciBwcmF5ZXIKaDM6IFNvY2lhbCBtZWRpYSBpcyBub3QgbWUKcDogRmFjZWJvb2ssIHR3aXR0ZXIsIGluc3RhZ3JhbSBhbmQgb3RoZXIgc29jaWFsIG1lZGlhIHBsYXRmb3JtcyBkbyBub3QgZ2l2ZSBvdXIgdHJ1ZSBwaWN0dXJlLiBUaGVzZSBwbGF0Zm9ybXMgZG8gbm90IHNob3cgdXMgZWF0aW5nIG1lYWxzLCB0YWxraW5nIHRvIG91ciBmYW1pbGllcywgZG9pbmcgb3VyIGpvYnMsIG1haW50YWluaW5nIHNvY2lhbCBjaXJjbGVzIGFuZCBkbyBleGVyY2lzZS4gV2Ugc2hhcmUgYSBwaG90byBvZiBzb21lIHBhcnQgb2Ygb3VyIGRheSwgbWFrZXMgcGVvcGxlIGp1ZGdlIG91ciBlbnRpcmUgbGlmZS4gV2UgZ28gb24gYSB0cmlwLCBzaGFyZSBpdHMgcGljdHVyZXMgd2l0aCBhIHNtaWxlIG9uIG91ciBmYWNlLCBkb2VzIG5vdCBtZWFuIHdlIGFyZSBhbHdheXMgc21pbGluZyBvbiBhIHRyaXAuPGI+IEl0IGlzIGEgbW9tZW50IGZyb20gdGhlIGVudGlyZSBkYXnigJlzIHRyaXAuIDwvYj4KcDogU2ltaWxhcmx5LCB3aGlsZSB3ZSBwb3N0IGEgcGljdHVyZSBvbiBzb2NpYWwgbWVkaWEsIHdlIHdhbnQgb3VyIGF1ZGllbmNlIChmcmllbmRzLCBmYW1pbHkgbWVtYmVycykgdG8gYmVsaWV2ZSB3ZSBhcmUgYWx3YXlzIHJvYW1pbmcgYXJvdW5kIHdpdGggYSBoYXBweSBmYWNlLiA8Yj5ObyBib2R5IGlzIGFsd2F5cyBoYXBweSwgZHJpdmluZyBoaXMgY2FyIHdpdGggaGlzIHdpZmU8L2I+LiBJdCBpcyBub3QgcG9zc2libGUuIFRoZSBtb21lbnQgd2Ugc2hvdyBvdXIgZW50aXJlIGRheSB0byBvdXIgZnJpZW5kcywgdGhleSB3aWxsIHJlYWxpc2UgaG93IHVnbHkgb3VyIGRheXMgZ28uIFRoZXNlIHBvc3RzIG9uIHNvY2lhbCBtZWRpYSBkb2VzIG5vdCBzaG93IHVzIGZpZ2h0aW5nLCBiZWluZyBzYWQsIHdlYXJpbmcgdGhlIGJhZCBjbG90aGVzLCBoYXZpbmcgYmFkIHRlZXRoIGFuZCBhIGJhZCBzbWVsbC4gSXQgb25seSBzaG93cyBvbmUgc2lkZSBvZiBhIGJpZ2dlciBwaWN0dXJlLgpwOiBJIGFtIHdyaXRpbmcgdGhpcyBibG9nIGJlZm9yZSBGYWphciBwcmF5ZXIuIEkgd3JpdGUgdGhpcyBibG9nIGluIDMwIHRvIDQwIG1pbnV0ZXMgd2hpbGUgdGhlIGRheSBoYXMgbm90IHlldCBzdGFydGVkLiBCeSB0aGUgdGltZSBJIGFtIGRvbmUgd3JpdGluZyBhbmQgcHVibGlzaGluZyBpdCwgSSBoYXZlIG5vdCByZWFsbHkgc3RhcnRlZCBteSBkYXkuIEkgaGF2ZSBub3QgeWV0IHByZXBhcmVkIG15c2VsZiBmb3Igb2ZmaWNlLCB0YWxrZWQgdG8gbXkgd2lmZSwgc29jaWFsaXNlZCB3aXRoIG15IGZyaWVuZHMgYW5kIHdhdGNoZWQgYSBtb3ZpZS4gSXQgaXMganVzdCB0aGUgNDAgbWludXRlcyBwYXJ0IG9mIG15IGVudGlyZSBkYXkuIApwOiBCdXQsIHRoZSBtb21lbnQsIEkgc2hhcmUgdGhpcyBibG9nIHdpdGggcGVvcGxlIG9uIGZhY2Vib29rLCB0d2l0dGVyIGFuZCBXaGF0c0FwcCBldmVyeWRheSwgc29tZSBwZW9wbGUgc3RhcnQgZ2V0dGluZyBqdWRnZW1lbnRhbC4gVGhleSBwZXJjZWl2ZSwg4oCcPGI+SSBhbSB3cml0aW5nIGl0IDI0IC8gNywgcm91bmQgdGhlIGNsb2NrPC9iPuKAnS4gSWYgSSBhbSBub3Qgd3JpdGluZyBpdCwg4oCcPGI+SSBhbSB0aGlua2luZyBhYm91dCBpdC4gSSBhbSBub3QgZ2V0dGluZyBmdWxsIHNsZWVwLiBJIGFtIGxpdmluZyBhIG1pc2VyYWJsZSBsaWZlPC9iPuKAnS4gQnV0LCBpbiByZWFsLCBJIGFtIGxpdmluZyBhIGxpZmUgb2YgZmlnaHRzLCB1cHNldHMsIGJvcmVkb21zIGFuZCBtb3ZpZXMuIApwOiBTb2NpYWwgbWVkaWEgZG8gbm90IHJlcHJlc2VudCBtZS4gSXQgZG9lc27igJl0IHJlcHJlc2VudCBhbnlib2R5LiBUaGUgZ3JlYXRlc3QgY2VsZWJyaXRpZXMgb24gdGhpcyBlYXJ0aCwgd2l0aCBhIHZlcmlmaWVkIGFjY291bnQsIGFyZSBvbmx5IHNob3dpbmcgb25lIHNpZGUgb2YgdGhlaXIgc3RvcnkuIFRoZWlyIGFjY291bnRzIGRvZXMgbm90IHJlcHJlc2VudCB0aGVtLiBXZSBhbGwgYXJlIGh1bWFuIGJlaW5ncy4gV2UgYWxsIGFyZSBqdXN0IGxpa2UgZWFjaCBvdGhlcjsgc2lubmluZywgc21pbGluZywgY3J5aW5nIGFuZCBmaWdodGluZy4gV2UgYXJlIG5vIGRpZmZlcmVudCBmcm9tIGVhY2ggb3RoZXIuCmgzOiBPbmUgc2lkZSBvZiB0aGUgcGljdHVyZQpwOiBMZXQgdXMgZ28gYSBsaXR0bGUgaGlnaGVyIGFuZCBzZWUgZnJvbSBhYm92ZS4gT3VyIGF0dGl0dWRlIG9mIDxiPmJlaW5nIGp1ZGdlbWVudGFsPC9iPiBoYXMgaW1wYWN0ZWQgb3VyIGxpdmVzLiBXZSBoYXZlIHRha2VuIGEgYm94LCBsYWJlbGxlZCBpdCBhbmQgcGxhY2VkIHRoaW5ncyBpbnNpZGUgdGhlc2UgYm94ZXMuCnA6IFNvbWUgYXJlIDxiPmV4dHJlbWlzdDwvYj4sIHNvbWUgYXJlPGI+IGJhZDwvYj4sIHNvbWUgYXJlPGI+IHN1Y2Nlc3NmdWw8L2I+IGFuZCBzb21lIGFyZTxiPiBmYWlsdXJlczwvYj4uIERpZmZlcmVudCBib3hlcywgZGlmZmVyZW50IG5hbWVzIGFuZCBkaWZmZXJlbnQgcGVvcGxlIGluIHRoZXNlIGJveGVzLiBGb3IgZXhhbXBsZQpoNTogMS4gUXVyYW4gaXMgYW4gZXh0cmVtaXN0IGJvb2sgb3IgUXVyYW4gaXMgYSBib29rIGRpZmZpY3VsdCB0byBhY3QgdXBvbgpwOiBXZSBoYXZlIG5vdCByZWFkIHRoZSBRdXJhbiBlbnRpcmVseS4gU29tZSBvZiB1cyBhcmUgbXVzbGltcywgYmVjYXVzZSB3ZSBhcmUgYXBwb2ludGVkIGEgPGI+TXVzbGltPC9iPiBieSBvdXIgcGFyZW50cy4gU29tZSBvZiB1cyBhcmUgY2hyaXN0aWFucywgYmVjYXVzZSB3ZSBhcmUgbGFiZWxsZWQgc28gYnkgb3VyIHBhcmVudHMuIApwOiBOb25lIG9mIHVzIGtub3c7IHdoYXQgUXVyYW4sIEJpYmxlLCBUb3JhaCBhbmQgWmFidXIgc3BlYWsgYWJvdXQuIFdlIGhhdmUgcmVhZCBvbmx5IDEgdmVyc2lvbiAvIDEgdHJhbnNsYXRpb24gb2YgdGhlc2UgYm9va3MuIFdlIGhhdmUgbGFiZWxsZWQgaXQg4oCYYm9yaW5n4oCZLCDigJhleHRyZW1pc3TigJksIOKAmGRpZmZpY3VsdCB0byBhY3QgdXBvbuKAmSBib29rLgpwOiBObyBib2R5IGtub3dzLCBub3QgZXZlbiB0aGUgZ3JlYXRlc3Qgc2Nob2xhciBrbm93LCB3aGF0IElzbGFtIG1lYW5zLiA8Yj5Jc2xhbSBoYXMgbm90IHlldCBtZXQgaXRzIGVuZC4gSXQgaXMgc3RpbGwgYmVpbmcgZGlzY292ZXJlZC48L2I+IApoNTogMi4gUGFraXN0YW4gaXMgYSBmYWlsaW5nIGVjb25vbXkKcDogV2UgYXJlIG5vdCBhbHdheXMgZmFpbGluZy4gT3VyIHN0b2NrIG1hcmtldCBnb2VzIHVwIHNvbWUgcGFydHMgb2YgdGhlIHllYXIuIEl0IGlzIGxvd2VyIHRvZGF5LCBidXQgNCBtb250aHMgYmFjayBpdCB3ZW50IHVwIHRvIDQwMDAwIHBvaW50cy4gSXQgZ29lcyB1cCwgY29tZXMgZG93biBhbmQgYWdhaW4gZ29lcyB1cC4gV2UgY2FuIG5vdCBmaXhhdGUgdGhlIGxhYmVsIG9mIGEg4oCYPGI+ZmFpbGluZzwvYj7igJkgZWNvbm9teSB0byBQYWtpc3Rhbi4gPGI+UGFraXN0YW4gaXMgYSBmbG91cmlzaGluZyBzdGF0ZSwgbmFtZWQgdGhlIGJlc3QgdG91cmlzdCBwbGFjZSBmb3IgeWVhciAyMDIwLjwvYj4KaDU6IDMuIFRoZSBib3hlciB3aG8gcHVuY2hlZCB0aGUgd2VhayBndXkgdGhhdCBoaXMgYmxvb2QgZ3VzaGVzIG91dApwLm5vdGU6IFRoaXMgcGFydCB3YXMgbmFycmF0ZWQgYnkgbXkgYnJvdGhlciB0aGF0IGhlIHNhdyBhIGJveGVyIGluIGEgY2x1YiwgcHVuY2hpbmcgYSBwb29yIGd1eS4gSGlzIGJsb29kIGd1c2hlZCBvcGVuIGFuZCBwZW9wbGUgZ290IG1hZCBhdCB0aGUgYm94ZXIuCnA6IFRoZSBib3hlciBpcyBub3QgYWx3YXlzIHB1bmNoaW5nIHRoZSB3ZWFrZXIgb25lcy4gSGUgcGF5cyBoaXMgZmFyZXMsIHRhbGtzIHRvIGhpcyBmYW1pbHksIGVhcm5zIG1vbmV5IGFuZCBzcGVuZHMgaXQgb24gaGltc2VsZi4gSGUgc2xlZXBzLCB0YWxrcywgY3JpZXMgYW5kIGxpdmVzIGEgbGlmZS4gSGlzIG9uZSB0aW1lIGhpdHRpbmcgYSBwZXJzb24sIGRvZXMgbm90IG1ha2UgaGltIGFuIGV2aWwgcGVyc29uLiBIZSB0dXJuZWQgZXZpbCBpbiBhIG1vbWVudCwgYW5kIGxhdGVyIGhlIGlzIGEgbm9ybWFsIGh1bWFuIGJlaW5nLiBIaXMgbW9zdCBhY3Rpb25zIGFyZSBsb3ZlYWJsZSBhbmQgZm9yZ2l2ZS1hYmxlLgpoNTogNC4gVGhlIHByb3N0aXR1dGUgd2lsbCBnbyB0byBoZWxsCnA6IE5vIHNoZSB3aWxsIG5vdCBnbyB0byBoZWxsLiBXZSBkbyBub3Qga25vdyB3aGF0IHRoZXkgZG8uIEhvdyB0aGV5IGxpdmUgdGhlaXIgbGl2ZXMuIFRoZXkgbWlnaHQgYmUgYXJlIG1vcmUgaHVtYmxlLCBtb3JlIGh1bWFuIGFuZCBtb3JlIGxvdmVhYmxlIHRoYW4gdXMuCmg1OiA1LiBBdGhlaXN0cyB3aWxsIGdvIHRvIGhlbGwKcDogV2UgZG8gbm90IGtub3cgeWV0LiBXZSBkbyBub3Qga25vdyB0aGUgUXVyYW4geWV0LiBXZSBkbyBub3Qga25vdyBob3cgQWxsYWggd29ya3MuIEhvdyBjYW4gd2UgYmUgc28gc3VyZSB0aGV5IHdpbGwgZ28gdG8gaGVsbD8gVGhleSBtaWdodCBnbyB0byBwYXJhZGlzZS4gIApoMzogU3RvcCBpdApwOiBXZSBuZWVkIHRvIHN0b3AgbGFiZWxsaW5nIHBlb3BsZS4gV2UgZG8gbm90IGtub3cgd2hhdCBhbnkgb2YgdXMgaXMgZ29pbmcgdGhyb3VnaC4gCnA6IEkgYW0gYSBzbWFsbCBraWQgaW5zaWRlIG15c2VsZi4gSSBrbm93IG15c2VsZiwgd2hvIG5lZWRzIHRvIGJlIGxvdmVkLCBjb25zaWRlcmVkIGh1bWFuIGJlaW5nIGFuZCBuZWVkcyB0byBiZSBlbmNvdXJhZ2VkLiBJIGtub3cgd2UgYWxsIGFyZSBsaWtlIG1lLCA8Yj5zbWFsbCBpbnNpZGU8L2I+LiBXZSBuZWVkIGhlbHAsIGdpdmUgdXMgaGVscC4gR2l2ZSB1cyBsb3ZlLCB0aGF0IHdlIGxpc3RlbiB0byB5b3UuIApwLmxhdmVuZGFyOiBNeSBibG9nIGlzIG9ubHkgNDAgbWludXRlcyBvZiBteSBlbnRpcmUgZGF5LiBJIGFtIG5vdCB0aGVzZSBibG9ncy4gCnAubm90ZTogR3V5cywgaWYgeW91IGxpa2UgdGhpcyBwcm9qZWN0LiBQbGVhc2UgZm9sbG93IHRoaXMgcHJvamVjdCdzIHBhZ2Ugb24gdHdpdHRlci4gPGEgaHJlZj0iaHR0cHM6Ly90d2l0dGVyLmNvbS96YWthdGxpc3RzIj5DbGljayBoZXJlIHRvIGdvIHRvIHRoZSB0d2l0dGVyIHBhZ2U8L2E+LiBTdWJzY3JpYmUgYmVsb3cgdG8gcmVjZWl2ZSB1cGRhdGVzLg=="
$bodyBytes = [System.Convert]::FromBase64String($bodyB64)
$body = [System.Text.Encoding]::UTF8.GetString($bodyBytes)

$row = 49

# Order of first-write matters for shared-string allocation order
# (Ser, Date set first since they are plain numbers; Tags text before the
# body text so the new shared strings land as 135, 136, 137 respectively).
$ws.Range("A$row").Value = 48
$ws.Range("B$row").Value = 43878
$ws.Range("B$row").NumberFormat = "d-mmm-yy"
$ws.Range("C$row").Value = "Surah Al Nisa, 51 - 67"
$ws.Range("F$row").Value = "Social Media, True picture of ourself, Athiests, Prostitutes, Pakistan, Quran and the holy books"
$ws.Range("D$row").Value = $body
$ws.Range("E$row").Value = "Qasim Ali"

$ws.Rows.Item($row).RowHeight = 409.6

$ws.Range("D$row").Select() | Out-Null

Write-Host "Row $row written, body length $($body.Length)"
